$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-22 down to 10-23
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new bibliography entry
$ws.Range("A9").Value2 = "Cerda R"
$ws.Range("B9").Value2 = "Silva A"
$ws.Range("C9").Value2 = "Valente J"
$ws.Range("F9").Value2 = "Impact of economic uncertainty in a small open economy: the case of Chile"
$ws.Range("G9").Value2 = "Paper"
$ws.Range("H9").Value2 = "Applied Economics"
$ws.Range("I9").Value2 = "Macroeconomía"
$ws.Range("E9").Value2 = 2018

# Rebuild the hyperlinks collection so refs line up with the shifted rows
# (the runtime does not auto-shift stored hyperlink ranges on row insert)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("J2"), "https://ideas.repec.org/p/pra/mprapa/79809.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J3"), "https://ideas.repec.org/p/pra/mprapa/83154.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J4"), "https://doi.org/10.1080/02692171.2019.1645816") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J5"), "https://www.bcentral.cl/en/content/-/details/monetary-policy-report-june-2015") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J6"), "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4043") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J7"), "https://www.bcentral.cl/en/web/banco-central/content/-/detalle/documento-de-trabajo-n-883") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J8"), "https://ideas.repec.org/a/chb/bcchni/v15y2012i1p105-117.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J9"), "https://ideas.repec.org/a/taf/applec/v50y2018i26p2894-2908.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J10"), "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4042") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J11"), "https://www.bcentral.cl/en/web/banco-central/content/-/detalle/analisis-de-sentimiento-basado-en-el-informe-de-percepciones-de-negocios-del-banco-central-de-chile") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J12"), "https://www.bcentral.cl/en/content/-/details/working-papers-n-888") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J13"), "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4041") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J14"), "https://ideas.repec.org/p/chb/bcchep/56.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J15"), "https://www.bcentral.cl/en/content/-/details/working-papers-n-876") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J16"), "https://www.bcentral.cl/en/content/-/details/working-papers-n-825") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J17"), "https://www.bcentral.cl/en/content/-/details/working-papers-n-899") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J18"), "https://ideas.repec.org/a/cml/moneta/vxxxiiy2009i2p181-208.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J19"), "https://www.google.com/url?sa=t&rct=j&q=&esrc=s&source=web&cd=&ved=2ahUKEwiE7Y60i5TvAhVQErkGHfecC4MQFjABegQIAhAD&url=http%3A%2F%2Fwww.oecd.org%2Fstd%2Fleading-indicators%2F43815334.pdf&usg=AOvVaw3BstLuhLtAOtjJeL5SsMj4") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J20"), "https://www.sciencedirect.com/science/article/abs/pii/S0169207019300676") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J21"), "https://www.bcentral.cl/en/content/-/details/working-papers-n-889") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J22"), "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/3564") | Out-Null
$ws.Hyperlinks.Add($ws.Range("J23"), "https://repositoriodigital.bcentral.cl/xmlui/handle/20.500.12580/4867") | Out-Null

# Restore hyperlink-style formatting (font/underline/color) on column J
$ws.Range("J2:J23").Style = "Hipervínculo"

# Refresh the recorded sort state to cover the extended data range A2:J23
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A23")) | Out-Null
$sortObj.SortFields.Add($ws.Range("E2:E23")) | Out-Null
$sortObj.SetRange($ws.Range("A2:J23"))
$sortObj.Header = -4142
$sortObj.Apply()
